$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "62.470.95"
Set-TextValue $ws.Range("E2") "  +2.20%  "

Set-TextValue $ws.Range("D3") "2.942.90"
Set-TextValue $ws.Range("E3") "  +1.91%  "

Set-TextValue $ws.Range("E4") "  -0.03%  "

Set-TextValue $ws.Range("D5") "589.16"
Set-TextValue $ws.Range("E5") "  +0.14%  "

Set-TextValue $ws.Range("D6") "146.83"
Set-TextValue $ws.Range("E6") "  +5.39%  "

Set-TextValue $ws.Range("E7") "  -0.03%  "

Set-TextValue $ws.Range("D8") "2.943.11"
Set-TextValue $ws.Range("E8") "  +1.94%  "

Set-TextValue $ws.Range("E9") "  +2.87%  "

Set-TextValue $ws.Range("D10") "7.07"
Set-TextValue $ws.Range("E10") "  +3.48%  "

Set-TextValue $ws.Range("D11") "0.151"
Set-TextValue $ws.Range("E11") "  +9.45%  "

Set-TextValue $ws.Range("E12") "  +1.49%  "

Set-TextValue $ws.Range("E13") "  +7.07%  "

Set-TextValue $ws.Range("D14") "32.24"
Set-TextValue $ws.Range("E14") "  -0.06%  "

Set-TextValue $ws.Range("D15") "0.126"
Set-TextValue $ws.Range("E15") "  -0.71%  "

Set-TextValue $ws.Range("D16") "3.428.84"
Set-TextValue $ws.Range("E16") "  +1.91%  "

Set-TextValue $ws.Range("D17") "62.434.49"
Set-TextValue $ws.Range("E17") "  +2.30%  "

Set-TextValue $ws.Range("D18") "6.66"
Set-TextValue $ws.Range("E18") "  +2.51%  "

Set-TextValue $ws.Range("D19") "2.936.67"
Set-TextValue $ws.Range("E19") "  +1.79%  "

Set-TextValue $ws.Range("D20") "433.87"
Set-TextValue $ws.Range("E20") "  +2.13%  "

Set-TextValue $ws.Range("D21") "13.47"
Set-TextValue $ws.Range("E21") "  +1.31%  "

Set-TextValue $ws.Range("D22") "0.663"
Set-TextValue $ws.Range("E22") "  +1.57%  "

Set-TextValue $ws.Range("D23") "6.97"
Set-TextValue $ws.Range("E23") "  +0.90%  "

Set-TextValue $ws.Range("D24") "11.14"
Set-TextValue $ws.Range("E24") "  +5.92%  "

Set-TextValue $ws.Range("D25") "80.12"

Set-TextValue $ws.Range("E26") "  +4.63%  "

Set-TextValue $ws.Range("D27") "2.10"
Set-TextValue $ws.Range("E27") "  +2.48%  "

Set-TextValue $ws.Range("E28") "  +0.02%  "

Set-TextValue $ws.Range("D29") "7.17"
Set-TextValue $ws.Range("E29") "  +8.40%  "

Set-TextValue $ws.Range("B30") "ImmutableX"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D30") "2.15"
Set-TextValue $ws.Range("E30") "  +4.69%  "

Set-TextValue $ws.Range("B31") "PancakeSwap"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "2.57"
Set-TextValue $ws.Range("E31") "  +1.73%  "

Set-TextValue $ws.Range("D32") "0.0000100"
Set-TextValue $ws.Range("E32") "  +18.31%  "

Set-TextValue $ws.Range("D33") "0.109"
Set-TextValue $ws.Range("E33") "  +4.26%  "

Set-TextValue $ws.Range("D34") "26.18"
Set-TextValue $ws.Range("E34") "  +2.30%  "

Set-TextValue $ws.Range("E35") "  -0.01%  "

Set-TextValue $ws.Range("D36") "0.991"

Set-TextValue $ws.Range("D37") "5.58"
Set-TextValue $ws.Range("E37") "  +2.85%  "

Set-TextValue $ws.Range("D38") "3.00"
Set-TextValue $ws.Range("E38") "  +7.17%  "

Set-TextValue $ws.Range("D39") "49.66"
Set-TextValue $ws.Range("E39") "  +1.41%  "

Set-TextValue $ws.Range("D40") "2.01"
Set-TextValue $ws.Range("E40") "  +5.69%  "

Set-TextValue $ws.Range("D41") "8.36"
Set-TextValue $ws.Range("E41") "  +0.44%  "

Set-TextValue $ws.Range("D42") "0.115"
Set-TextValue $ws.Range("E42") "  -0.94%  "

Set-TextValue $ws.Range("D43") "0.275"
Set-TextValue $ws.Range("E43") "  +4.47%  "

Set-TextValue $ws.Range("D44") "38.57"
Set-TextValue $ws.Range("E44") "  +0.48%  "

Set-TextValue $ws.Range("D45") "135.09"
Set-TextValue $ws.Range("E45") "  +2.00%  "

Set-TextValue $ws.Range("D46") "2.690.73"
Set-TextValue $ws.Range("E46") "  +1.32%  "

Set-TextValue $ws.Range("D47") "0.0337"
Set-TextValue $ws.Range("E47") "  +2.56%  "

Set-TextValue $ws.Range("D48") "354.03"
Set-TextValue $ws.Range("E48") "  +1.74%  "

Set-TextValue $ws.Range("E49") "  +0.03%  "

Set-TextValue $ws.Range("E50") "  +1.96%  "

Set-TextValue $ws.Range("D51") "22.51"
Set-TextValue $ws.Range("E51") "  +1.15%  "
